# Kilimanjaro Weekly Scoreboard — append this week's workout rows (149-152)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formatting from the last existing data row (148) down into
# the four new rows so the Date column (B) keeps its date style, exactly
# like every other row in the table.
$ws.Range("A148:M148").Copy()
$ws.Range("A149:M152").PasteSpecial(-4122)

# Row 149 - Steven / Run / Agile Antelope week
$ws.Range("A149").Value = "Steven"
$ws.Range("B149").Value = 45475
$ws.Range("C149").Value = "Workout"
$ws.Range("D149").Value = 44
$ws.Range("E149").Value = 0
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 23
$ws.Range("H149").Value = 1
$ws.Range("I149").Value = 1
$ws.Range("J149").Value = 0
$ws.Range("K149").Value = 0
$ws.Range("L149").Value = "Mighty Monkey"
$ws.Range("M149").Value = 4

# Row 150 - Matt / Walk / Sauntering Hippo week
$ws.Range("A150").Value = "Matt"
$ws.Range("B150").Value = 45476
$ws.Range("C150").Value = "Walk"
$ws.Range("D150").Value = 21
$ws.Range("E150").Value = 0.97
$ws.Range("F150").Value = 56
$ws.Range("G150").Value = 21
$ws.Range("H150").Value = 0
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 0
$ws.Range("K150").Value = 0
$ws.Range("L150").Value = "Sauntering Hippo"
$ws.Range("M150").Value = 4

# Row 151 - Jeremiah / Workout / Agile Antelope week
$ws.Range("A151").Value = "Jeremiah"
$ws.Range("B151").Value = 45476
$ws.Range("C151").Value = "Workout"
$ws.Range("D151").Value = 56
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 54
$ws.Range("H151").Value = 2
$ws.Range("I151").Value = 0
$ws.Range("J151").Value = 0
$ws.Range("K151").Value = 0
$ws.Range("L151").Value = "Agile Antelope"
$ws.Range("M151").Value = 4

# Row 152 - Eric / Ride / Agile Antelope week
$ws.Range("A152").Value = "Eric"
$ws.Range("B152").Value = 45476
$ws.Range("C152").Value = "Run"
$ws.Range("D152").Value = 30
$ws.Range("E152").Value = 3.02
$ws.Range("F152").Value = 62
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1
$ws.Range("I152").Value = 3
$ws.Range("J152").Value = 16
$ws.Range("K152").Value = 7
$ws.Range("L152").Value = "Agile Antelope"
$ws.Range("M152").Value = 4

# Match the saved selection state: active cell moves to E153 (just below the
# newly-added data), on the frozen bottom-left pane.
[void]$ws.Range("E153").Select()
